$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the A (code) and D (level) values for the new rows first.
$ws.Cells.Item(45, 1).Value = 10043
$ws.Cells.Item(46, 1).Value = 10044
$ws.Cells.Item(47, 1).Value = 10045
$ws.Cells.Item(48, 1).Value = 10046
$ws.Cells.Item(49, 1).Value = 10047
$ws.Cells.Item(50, 1).Value = 10048

# Set the B (message) values in the exact order the new shared strings were
# originally introduced, so the resulting shared-string table order matches.
$ws.Cells.Item(45, 2).Value = "message_10043_phone_address_record_created_successfully"
$ws.Cells.Item(47, 2).Value = "message_10045_phone_address_record_deleted_successfully"
$ws.Cells.Item(46, 2).Value = "message_10044_phone_address_record_updated_successfully"
$ws.Cells.Item(48, 2).Value = "message_10046_electronic_address_record_created_successfully"
$ws.Cells.Item(49, 2).Value = "message_10047_electronic_address_record_updated_successfully"
$ws.Cells.Item(50, 2).Value = "message_10048_electronic_address_record_deleted_successfully"

$ws.Cells.Item(45, 4).Value = "Success"
$ws.Cells.Item(46, 4).Value = "Success"
$ws.Cells.Item(47, 4).Value = "Success"
$ws.Cells.Item(48, 4).Value = "Success"
$ws.Cells.Item(49, 4).Value = "Success"
$ws.Cells.Item(50, 4).Value = "Success"

# Update visible window / selection to match end-state view
$ws.Application.ActiveWindow.ScrollRow = 42
$ws.Range("B50").Select()
